# Updating the GESS model with the latest data
# Shift column A timestamps by +5 days (46074.x -> 46079.x, 46075 -> 46080)
# and replace column B values with the latest Actual Production (MW) data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newA = @(46079.01041666666,46079.02083333334,46079.03125,46079.04166666666,46079.05208333334,46079.0625,46079.07291666666,46079.08333333334,46079.09375,46079.10416666666,46079.11458333334,46079.125,46079.13541666666,46079.14583333334,46079.15625,46079.16666666666,46079.17708333334,46079.1875,46079.19791666666,46079.20833333334,46079.21875,46079.22916666666,46079.23958333334,46079.25,46079.26041666666,46079.27083333334,46079.28125,46079.29166666666,46079.30208333334,46079.3125,46079.32291666666,46079.33333333334,46079.34375,46079.35416666666,46079.36458333334,46079.375,46079.38541666666,46079.39583333334,46079.40625,46079.41666666666,46079.42708333334,46079.4375,46079.44791666666,46079.45833333334,46079.46875,46079.47916666666,46079.48958333334,46079.5,46079.51041666666,46079.52083333334,46079.53125,46079.54166666666,46079.55208333334,46079.5625,46079.57291666666,46079.58333333334,46079.59375,46079.60416666666,46079.61458333334,46079.625,46079.63541666666,46079.64583333334,46079.65625,46079.66666666666,46079.67708333334,46079.6875,46079.69791666666,46079.70833333334,46079.71875,46079.72916666666,46079.73958333334,46079.75,46079.76041666666,46079.77083333334,46079.78125,46079.79166666666,46079.80208333334,46079.8125,46079.82291666666,46079.83333333334,46079.84375,46079.85416666666,46079.86458333334,46079.875,46079.88541666666,46079.89583333334,46079.90625,46079.91666666666,46079.92708333334,46079.9375,46079.94791666666,46079.95833333334,46079.96875,46079.97916666666,46079.98958333334,46080)

$newB = @(856,853,866,0,833,812,807,816,863,882,891,924,922,974,1024,1084,1135,1183,1181,0,1266,1309,1356,1391,1432,1444,1383,1355,1365,1437,1531,1566,1545,1595,1626,1552,1439,1388,1401,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0)

for ($i = 0; $i -lt $newB.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value2 = $newA[$i]
    $ws.Cells.Item($row, 2).Value2 = $newB[$i]
}
